$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range('D2').Value = 'maa://24702 (94.27), maa://25390 (95.92), maa://36681 (87.01)'
$ws.Range('AF2').Value = 'maa://25251 (91.67), ***maa://21730 (20.59), ***maa://39501 (18.18), *maa://36675 (60.0)'
$ws.Range('L3').Value = '*maa://22880 (65.76), maa://20276 (84.91), *maa://22749 (72.73)'
$ws.Range('T3').Value = 'maa://24617 (88.99), **maa://20790 (44.78), ***maa://37170 (19.3)'
$ws.Range('X3').Value = 'maa://27396 (84.42), maa://27484 (96.15), maa://27480 (82.86)'
$ws.Range('D4').Value = 'maa://24632 (93.29), **maa://24303 (33.33), maa://22499 (86.67), maa://22746 (100.0)'
$ws.Range('AB5').Value = '*maa://29863 (69.7), ***maa://22752 (13.33), **maa://26013 (42.86)'
$ws.Range('A8').Value = '更新日期：2024.12.24 13:17:51'
$ws.Range('X8').Value = 'maa://21411 (95.92)'
$ws.Range('X9').Value = 'maa://26223 (97.48)'
$ws.Range('P10').Value = 'maa://28977 (91.36), maa://36669 (89.47), *maa://23264 (61.82)'
$ws.Range('T10').Value = 'maa://27395 (96.0), maa://22755 (87.61), **maa://22756 (40.91), ***maa://21737 (10.61)'
$ws.Range('X10').Value = 'maa://22301 (97.62), maa://22726 (100.0)'
$ws.Range('D15').Value = '*maa://22743 (76.96), maa://22734 (84.03), *maa://30808 (65.08), **maa://36048 (33.33)'
$ws.Range('T16').Value = 'maa://22729 (95.39), *maa://28648 (68.33), maa://36674 (82.93)'
$ws.Range('H17').Value = 'maa://22430 (88.65), maa://39599 (86.84)'
$ws.Range('H18').Value = 'maa://24421 (90.38)'
$ws.Range('T19').Value = 'maa://24386 (99.05)'
$ws.Range('D20').Value = 'maa://21432 (89.86), maa://25198 (93.07), *maa://20795 (51.18), maa://36680 (96.67)'
$ws.Range('L23').Value = 'maa://39756 (94.12), maa://39875 (93.75)'
$ws.Range('P23').Value = 'maa://30587 (92.02), *maa://29748 (75.78), ***maa://29785 (16.42), *maa://37566 (71.43)'
$ws.Range('X24').Value = 'maa://29988 (86.73), maa://23504 (93.28), **maa://22892 (39.58), *maa://25141 (76.98), maa://36663 (80.6), ***maa://22815 (23.08)'
$ws.Range('AF24').Value = 'maa://22523 (85.64), maa://36672 (80.77), maa://29910 (92.45), **maa://21440 (34.55)'
$ws.Range('X28').Value = 'maa://39929 (89.84), ***maa://39723 (14.29), maa://41749 (90.91)'
$ws.Range('AF28').Value = 'maa://36660 (92.93), *maa://36701 (64.29)'
$ws.Range('AB30').Value = 'maa://42979 (96.3)'
$ws.Range('L31').Value = 'maa://35926 (93.98), maa://36258 (83.16), *maa://43904 (75.0)'
$ws.Range('T32').Value = 'maa://42859 (96.2), maa://41108 (87.76), maa://41238 (96.15)'
$ws.Range('H39').Value = 'maa://25199 (84.82), maa://36670 (87.21), maa://30434 (89.23), ***maa://25036 (16.0), *maa://44165 (66.67)'
$ws.Range('H46').Value = 'maa://35931 (92.5), maa://43901 (88.89)'
